$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows before row 1051, shifting existing rows 1051-1142 down to 1053-1144
$ws.Rows("1051:1052").Insert()

# New row 1051 - "Primera" quality, date 45106
$ws.Cells.Item(1051, 1).Value = 3
$ws.Cells.Item(1051, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(1051, 3).Value = "Coquimbo"
$ws.Cells.Item(1051, 4).Value = 45106
$ws.Cells.Item(1051, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1051, 5).Value = 5
$ws.Cells.Item(1051, 6).Value = 100112023
$ws.Cells.Item(1051, 7).Value = "Brócoli"
$ws.Cells.Item(1051, 8).Value = "Sin especificar"
$ws.Cells.Item(1051, 9).Value = "Primera"
$ws.Cells.Item(1051, 10).Value = 3700
$ws.Cells.Item(1051, 11).Value = 750
$ws.Cells.Item(1051, 12).Value = 800
$ws.Cells.Item(1051, 13).Value = 776
$ws.Cells.Item(1051, 14).Value = "$/unidad"
$ws.Cells.Item(1051, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(1051, 16).Value = 776
$ws.Cells.Item(1051, 17).Value = 1
$ws.Cells.Item(1051, 18).Value = "Hortaliza"

# New row 1052 - "Segunda" quality, date 45106
$ws.Cells.Item(1052, 1).Value = 3
$ws.Cells.Item(1052, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(1052, 3).Value = "Coquimbo"
$ws.Cells.Item(1052, 4).Value = 45106
$ws.Cells.Item(1052, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(1052, 5).Value = 5
$ws.Cells.Item(1052, 6).Value = 100112023
$ws.Cells.Item(1052, 7).Value = "Brócoli"
$ws.Cells.Item(1052, 8).Value = "Sin especificar"
$ws.Cells.Item(1052, 9).Value = "Segunda"
$ws.Cells.Item(1052, 10).Value = 1800
$ws.Cells.Item(1052, 11).Value = 600
$ws.Cells.Item(1052, 12).Value = 600
$ws.Cells.Item(1052, 13).Value = 600
$ws.Cells.Item(1052, 14).Value = "$/unidad"
$ws.Cells.Item(1052, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(1052, 16).Value = 600
$ws.Cells.Item(1052, 17).Value = 1
$ws.Cells.Item(1052, 18).Value = "Hortaliza"
